# Avance en DAG para validadores
# Se crean procedimientos de validacion para rdc01 y rdc20,
# se empieza a trabajar en los DAG del validador

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Nuevo validador 14: Tipo de registro no corresponde (RDC20) ---
$ws.Range("A15").Value = "14"
$ws.Range("B15").Value = "Tipo de registro no corresponde"
$ws.Range("C15").Value = "Generico"
$ws.Range("D15").Value = "No aplica"
$ws.Range("E15").Value = "proceso.val_en_dominio"

# --- Nuevo validador 15: Tipo de flujo informado no corresponde ---
$ws.Range("A16").Value = "15"
$ws.Range("B16").Value = "Tipo de flujo informado no corresponde"
$ws.Range("C16").Value = "Puntual"
$ws.Range("D16").Value = "RDC20"
$ws.Range("E16").Value = "proceso.val_codigo_tabla"

# Columna Explicacion para los dos validadores recien creados
$ws.Range("F15").Value = "Corresponde a validar si los tipos de registros corresponden a lo solicitado"
$ws.Range("F16").Value = "Corresponde a validar el dato que solicita RDC20 para el campo tipo de flujo"

# --- Nuevo validador 16: Filler ---
$ws.Range("A17").Value = "16"
$ws.Range("E17").Value = "proceso.val_num_16"
$ws.Range("B17").Value = "Filler debe ser completado con espacios y tener un largo especifico"
$ws.Range("C17").Value = "Generico"
$ws.Range("D17").Value = "No aplica"
$ws.Range("F17").Value = "Corresponde a validar el filler si esta compuesto por espacios y tiene un largo especifico"

# Deja la seleccion activa donde el autor la dejo al guardar
$ws.Range("E7").Select()
